$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value2 = "ECs"
$ws.Cells.Item(2,2).Value2 = "Pomc"
$ws.Cells.Item(2,3).Value2 = "Mc4r"
$ws.Cells.Item(2,4).Value2 = "ECs"
$ws.Cells.Item(2,5).Value2 = 3
$ws.Cells.Item(2,6).Value2 = 1
$ws.Cells.Item(2,7).Value2 = 2.862146333333333
$ws.Cells.Item(2,8).Value2 = 8.586439
$ws.Cells.Item(2,9).Value2 = 0.9778268096017091
$ws.Cells.Item(2,10).Value2 = 0.9778268096017091
$ws.Cells.Item(2,11).Value2 = 1
$ws.Cells.Item(2,12).Value2 = 0.3333333333333333
$ws.Cells.Item(2,13).Value2 = 0.01244633333333333
$ws.Cells.Item(2,14).Value2 = 0.037339
$ws.Cells.Item(2,15).Value2 = 0.03341097222806328
$ws.Cells.Item(2,16).Value2 = 0.03341097222806328
$ws.Cells.Item(2,17).Value2 = 0.03562322731344444
$ws.Cells.Item(2,18).Value2 = 0.320609045821
$ws.Cells.Item(2,19).Value2 = 0.03267014437945842
$ws.Cells.Item(2,20).Value2 = 0.03267014437945842
$ws.Cells.Item(3,1).Value2 = "ECs"
$ws.Cells.Item(3,2).Value2 = "Pomc"
$ws.Cells.Item(3,3).Value2 = "Mc4r"
$ws.Cells.Item(3,4).Value2 = "FAPs"
$ws.Cells.Item(3,5).Value2 = 3
$ws.Cells.Item(3,6).Value2 = 1
$ws.Cells.Item(3,7).Value2 = 2.862146333333333
$ws.Cells.Item(3,8).Value2 = 8.586439
$ws.Cells.Item(3,9).Value2 = 0.9778268096017091
$ws.Cells.Item(3,10).Value2 = 0.9778268096017091
$ws.Cells.Item(3,11).Value2 = 2
$ws.Cells.Item(3,12).Value2 = 0.6666666666666666
$ws.Cells.Item(3,13).Value2 = 0.05317933333333334
$ws.Cells.Item(3,14).Value2 = 0.159538
$ws.Cells.Item(3,15).Value2 = 0.1427547520640821
$ws.Cells.Item(3,16).Value2 = 0.1427547520640821
$ws.Cells.Item(3,17).Value2 = 0.1522070339091111
$ws.Cells.Item(3,18).Value2 = 1.369863305182
$ws.Cells.Item(3,19).Value2 = 0.1395894237663043
$ws.Cells.Item(3,20).Value2 = 0.1395894237663044
$ws.Cells.Item(4,1).Value2 = "ECs"
$ws.Cells.Item(4,2).Value2 = "Pomc"
$ws.Cells.Item(4,3).Value2 = "Mc4r"
$ws.Cells.Item(4,4).Value2 = "MuSCs"
$ws.Cells.Item(4,5).Value2 = 3
$ws.Cells.Item(4,6).Value2 = 1
$ws.Cells.Item(4,7).Value2 = 2.862146333333333
$ws.Cells.Item(4,8).Value2 = 8.586439
$ws.Cells.Item(4,9).Value2 = 0.9778268096017091
$ws.Cells.Item(4,10).Value2 = 0.9778268096017091
$ws.Cells.Item(4,11).Value2 = 3
$ws.Cells.Item(4,12).Value2 = 1
$ws.Cells.Item(4,13).Value2 = 0.282443
$ws.Cells.Item(4,14).Value2 = 0.847329
$ws.Cells.Item(4,15).Value2 = 0.758190784087218
$ws.Cells.Item(4,16).Value2 = 0.7581907840872181
$ws.Cells.Item(4,17).Value2 = 0.8083931968256667
$ws.Cells.Item(4,18).Value2 = 7.275538771431
$ws.Cells.Item(4,19).Value2 = 0.7413792754734226
$ws.Cells.Item(4,20).Value2 = 0.7413792754734227
$ws.Cells.Item(5,1).Value2 = "ECs"
$ws.Cells.Item(5,2).Value2 = "Pomc"
$ws.Cells.Item(5,3).Value2 = "Mc4r"
$ws.Cells.Item(5,4).Value2 = "Resolving-Mac"
$ws.Cells.Item(5,5).Value2 = 3
$ws.Cells.Item(5,6).Value2 = 1
$ws.Cells.Item(5,7).Value2 = 2.862146333333333
$ws.Cells.Item(5,8).Value2 = 8.586439
$ws.Cells.Item(5,9).Value2 = 0.9778268096017091
$ws.Cells.Item(5,10).Value2 = 0.9778268096017091
$ws.Cells.Item(5,11).Value2 = 1
$ws.Cells.Item(5,12).Value2 = 0.3333333333333333
$ws.Cells.Item(5,13).Value2 = 0.02445366666666667
$ws.Cells.Item(5,14).Value2 = 0.073361
$ws.Cells.Item(5,15).Value2 = 0.06564349162063661
$ws.Cells.Item(5,16).Value2 = 0.06564349162063661
$ws.Cells.Item(5,17).Value2 = 0.06998997238655555
$ws.Cells.Item(5,18).Value2 = 0.629909751479
$ws.Cells.Item(5,19).Value2 = 0.06418796598252362
$ws.Cells.Item(5,20).Value2 = 0.06418796598252362
$ws.Cells.Item(6,1).Value2 = "MuSCs"
$ws.Cells.Item(6,2).Value2 = "Pomc"
$ws.Cells.Item(6,3).Value2 = "Mc4r"
$ws.Cells.Item(6,4).Value2 = "ECs"
$ws.Cells.Item(6,5).Value2 = 1
$ws.Cells.Item(6,6).Value2 = 0.3333333333333333
$ws.Cells.Item(6,7).Value2 = 0.064902
$ws.Cells.Item(6,8).Value2 = 0.194706
$ws.Cells.Item(6,9).Value2 = 0.02217319039829088
$ws.Cells.Item(6,10).Value2 = 0.02217319039829088
$ws.Cells.Item(6,11).Value2 = 1
$ws.Cells.Item(6,12).Value2 = 0.3333333333333333
$ws.Cells.Item(6,13).Value2 = 0.01244633333333333
$ws.Cells.Item(6,14).Value2 = 0.037339
$ws.Cells.Item(6,15).Value2 = 0.03341097222806328
$ws.Cells.Item(6,16).Value2 = 0.03341097222806328
$ws.Cells.Item(6,17).Value2 = 0.000807791926
$ws.Cells.Item(6,18).Value2 = 0.007270127333999999
$ws.Cells.Item(6,19).Value2 = 0.0007408278486048561
$ws.Cells.Item(6,20).Value2 = 0.0007408278486048561
$ws.Cells.Item(7,1).Value2 = "MuSCs"
$ws.Cells.Item(7,2).Value2 = "Pomc"
$ws.Cells.Item(7,3).Value2 = "Mc4r"
$ws.Cells.Item(7,4).Value2 = "FAPs"
$ws.Cells.Item(7,5).Value2 = 1
$ws.Cells.Item(7,6).Value2 = 0.3333333333333333
$ws.Cells.Item(7,7).Value2 = 0.064902
$ws.Cells.Item(7,8).Value2 = 0.194706
$ws.Cells.Item(7,9).Value2 = 0.02217319039829088
$ws.Cells.Item(7,10).Value2 = 0.02217319039829088
$ws.Cells.Item(7,11).Value2 = 2
$ws.Cells.Item(7,12).Value2 = 0.6666666666666666
$ws.Cells.Item(7,13).Value2 = 0.05317933333333334
$ws.Cells.Item(7,14).Value2 = 0.159538
$ws.Cells.Item(7,15).Value2 = 0.1427547520640821
$ws.Cells.Item(7,16).Value2 = 0.1427547520640821
$ws.Cells.Item(7,17).Value2 = 0.003451445092
$ws.Cells.Item(7,18).Value2 = 0.031063005828
$ws.Cells.Item(7,19).Value2 = 0.0031653282977777
$ws.Cells.Item(7,20).Value2 = 0.0031653282977777
$ws.Cells.Item(8,1).Value2 = "MuSCs"
$ws.Cells.Item(8,2).Value2 = "Pomc"
$ws.Cells.Item(8,3).Value2 = "Mc4r"
$ws.Cells.Item(8,4).Value2 = "MuSCs"
$ws.Cells.Item(8,5).Value2 = 1
$ws.Cells.Item(8,6).Value2 = 0.3333333333333333
$ws.Cells.Item(8,7).Value2 = 0.064902
$ws.Cells.Item(8,8).Value2 = 0.194706
$ws.Cells.Item(8,9).Value2 = 0.02217319039829088
$ws.Cells.Item(8,10).Value2 = 0.02217319039829088
$ws.Cells.Item(8,11).Value2 = 3
$ws.Cells.Item(8,12).Value2 = 1
$ws.Cells.Item(8,13).Value2 = 0.282443
$ws.Cells.Item(8,14).Value2 = 0.847329
$ws.Cells.Item(8,15).Value2 = 0.758190784087218
$ws.Cells.Item(8,16).Value2 = 0.7581907840872181
$ws.Cells.Item(8,17).Value2 = 0.018331115586
$ws.Cells.Item(8,18).Value2 = 0.164980040274
$ws.Cells.Item(8,19).Value2 = 0.01681150861379534
$ws.Cells.Item(8,20).Value2 = 0.01681150861379534
$ws.Cells.Item(9,1).Value2 = "MuSCs"
$ws.Cells.Item(9,2).Value2 = "Pomc"
$ws.Cells.Item(9,3).Value2 = "Mc4r"
$ws.Cells.Item(9,4).Value2 = "Resolving-Mac"
$ws.Cells.Item(9,5).Value2 = 1
$ws.Cells.Item(9,6).Value2 = 0.3333333333333333
$ws.Cells.Item(9,7).Value2 = 0.064902
$ws.Cells.Item(9,8).Value2 = 0.194706
$ws.Cells.Item(9,9).Value2 = 0.02217319039829088
$ws.Cells.Item(9,10).Value2 = 0.02217319039829088
$ws.Cells.Item(9,11).Value2 = 1
$ws.Cells.Item(9,12).Value2 = 0.3333333333333333
$ws.Cells.Item(9,13).Value2 = 0.02445366666666667
$ws.Cells.Item(9,14).Value2 = 0.073361
$ws.Cells.Item(9,15).Value2 = 0.06564349162063661
$ws.Cells.Item(9,16).Value2 = 0.06564349162063661
$ws.Cells.Item(9,17).Value2 = 0.001587091874
$ws.Cells.Item(9,18).Value2 = 0.014283826866
$ws.Cells.Item(9,19).Value2 = 0.001455525638112988
$ws.Cells.Item(9,20).Value2 = 0.001455525638112988
